$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 23 (Hero Canal Stop-Log Gate... row),
# shifting the existing rows 23-29 down to 24-30. By default the new row
# inherits row 22's formatting, so after inserting we copy the number
# format from the (now shifted) row 24 onto the new A23 cell to match the
# style used by the block of rows below it.
$ws.Rows(23).Insert()

$ws.Range("A24").Copy()
$ws.Range("A23").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New gate entry: Houma Navigation Canal
$ws.Range("A23").Value = "HoumaNavCanal"
$ws.Range("B23").Value = "Houma Navigation Canal"
$ws.Range("C23").Value = 0

# Leave the selection on the newly added description cell, matching the
# workbook author's last editing position.
$ws.Range("B23").Select() | Out-Null
